$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7845051884651184
$ws.Range("B1").Value = 3.319449186325073
$ws.Range("C1").Value = 3.656013250350952
$ws.Range("D1").Value = 2.975138902664185
$ws.Range("E1").Value = 1.783617615699768
